$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.718.85"
$ws.Cells.Item(2, 5).Value = "  +2.33%  "

$ws.Cells.Item(3, 4).Value = "1.871.32"
$ws.Cells.Item(3, 5).Value = "  +2.01%  "

$ws.Cells.Item(4, 5).Value = "  +0.28%  "

$ws.Cells.Item(5, 4).Value = "324.45"
$ws.Cells.Item(5, 5).Value = "  +0.05%  "

$ws.Cells.Item(6, 5).Value = "  +0.21%  "

$ws.Cells.Item(7, 4).Value = "0.4611"
$ws.Cells.Item(7, 5).Value = "  -0.59%  "

$ws.Cells.Item(8, 4).Value = "0.3857"
$ws.Cells.Item(8, 5).Value = "  -0.22%  "

$ws.Cells.Item(9, 4).Value = "0.07862"
$ws.Cells.Item(9, 5).Value = "  +0.01%  "

$ws.Cells.Item(10, 4).Value = "0.9774"
$ws.Cells.Item(10, 5).Value = "  +1.74%  "

$ws.Cells.Item(11, 4).Value = "21.81"
$ws.Cells.Item(11, 5).Value = "  -0.45%  "

$ws.Cells.Item(12, 4).Value = "1.861.00"
$ws.Cells.Item(12, 5).Value = "  +3.20%  "

$ws.Cells.Item(13, 4).Value = "7.002"
$ws.Cells.Item(13, 5).Value = "  +1.24%  "

$ws.Cells.Item(14, 4).Value = "5.702"
$ws.Cells.Item(14, 5).Value = "  +0.29%  "

$ws.Cells.Item(15, 4).Value = "0.06942"
$ws.Cells.Item(15, 5).Value = "  +1.42%  "

$ws.Cells.Item(16, 4).Value = "88.41"
$ws.Cells.Item(16, 5).Value = "  +1.07%  "

$ws.Cells.Item(17, 5).Value = "  +0.30%  "

$ws.Cells.Item(18, 4).Value = "0.00001002"
$ws.Cells.Item(18, 5).Value = "  +0.86%  "

$ws.Cells.Item(19, 4).Value = "16.77"
$ws.Cells.Item(19, 5).Value = "  +0.42%  "

$ws.Cells.Item(20, 4).Value = "1.002"
$ws.Cells.Item(20, 5).Value = "  +0.29%  "

$ws.Cells.Item(21, 4).Value = "28.705.21"
$ws.Cells.Item(21, 5).Value = "  +2.24%  "

$ws.Cells.Item(22, 4).Value = "5.269"
$ws.Cells.Item(22, 5).Value = "  -1.12%  "

$ws.Cells.Item(23, 4).Value = "11.07"
$ws.Cells.Item(23, 5).Value = "  +0.59%  "

$ws.Cells.Item(24, 4).Value = "2.096"
$ws.Cells.Item(24, 5).Value = "  +0.05%  "

$ws.Cells.Item(25, 4).Value = "2.041.21"
$ws.Cells.Item(25, 5).Value = "  +0.22%  "

$ws.Cells.Item(26, 4).Value = "152.73"
$ws.Cells.Item(26, 5).Value = "  -0.93%  "

$ws.Cells.Item(27, 4).Value = "19.30"
$ws.Cells.Item(27, 5).Value = "  +0.81%  "

$ws.Cells.Item(28, 4).Value = "5.869"
$ws.Cells.Item(28, 5).Value = "  +3.43%  "

$ws.Cells.Item(29, 4).Value = "1.988"
$ws.Cells.Item(29, 5).Value = "  +1.26%  "

$ws.Cells.Item(30, 4).Value = "119.10"
$ws.Cells.Item(30, 5).Value = "  +0.76%  "

$ws.Cells.Item(31, 4).Value = "0.09323"
$ws.Cells.Item(31, 5).Value = "  +0.97%  "

$ws.Cells.Item(32, 4).Value = "0.9167"
$ws.Cells.Item(32, 5).Value = "  -2.15%  "

$ws.Cells.Item(33, 4).Value = "5.293"
$ws.Cells.Item(33, 5).Value = "  +0.45%  "

$ws.Cells.Item(34, 4).Value = "1.332"
$ws.Cells.Item(34, 5).Value = "  +0.75%  "

$ws.Cells.Item(35, 5).Value = "  +0.55%  "

$ws.Cells.Item(36, 4).Value = "0.05789"
$ws.Cells.Item(36, 5).Value = "  -1.08%  "

$ws.Cells.Item(37, 5).Value = "  +0.99%  "

$ws.Cells.Item(38, 4).Value = "0.02078"
$ws.Cells.Item(38, 5).Value = "  -2.38%  "

$ws.Cells.Item(39, 4).Value = "7.647"
$ws.Cells.Item(39, 5).Value = "  -1.70%  "

$ws.Cells.Item(40, 4).Value = "0.5621"
$ws.Cells.Item(40, 5).Value = "  +0.58%  "

$ws.Cells.Item(41, 4).Value = "0.1781"
$ws.Cells.Item(41, 5).Value = "  +1.04%  "

$ws.Cells.Item(42, 4).Value = "9.767"
$ws.Cells.Item(42, 5).Value = "  -1.13%  "

$ws.Cells.Item(43, 4).Value = "0.07225"
$ws.Cells.Item(43, 5).Value = "  -0.75%  "

$ws.Cells.Item(44, 4).Value = "11.75"
$ws.Cells.Item(44, 5).Value = "  +0.55%  "

$ws.Cells.Item(45, 4).Value = "0.5281"
$ws.Cells.Item(45, 5).Value = "  +0.25%  "

$ws.Cells.Item(46, 4).Value = "2.131"
$ws.Cells.Item(46, 5).Value = "  +0.24%  "

$ws.Cells.Item(47, 4).Value = "1.121"
$ws.Cells.Item(47, 5).Value = "  +0.16%  "

$ws.Cells.Item(48, 4).Value = "1.835"
$ws.Cells.Item(48, 5).Value = "  +0.22%  "

$ws.Cells.Item(49, 4).Value = "112.83"
$ws.Cells.Item(49, 5).Value = "  +0.26%  "

$ws.Cells.Item(50, 4).Value = "2.411"
$ws.Cells.Item(50, 5).Value = "  +3.72%  "

$ws.Cells.Item(51, 4).Value = "1.002"
$ws.Cells.Item(51, 5).Value = "  +0.22%  "
